$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain text (e.g. "27.609.28", "  -1.29%  ").
# Pin the columns to Text format before writing so Excel does not
# auto-coerce numeric-looking strings (like "1.00" or "211.50") into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.542.02'
$ws.Range('E2').Value = '  -1.45%  '
$ws.Range('D3').Value = '1.631.02'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '211.50'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '23.05'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('D10').Value = '0.0610'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').Value = '0.0862'
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('D12').Value = '1.863.57'
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = '1.632.41'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '0.558'
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').Value = '65.03'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').Value = '27.508.68'
$ws.Range('D18').Value = '229.38'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').Value = '10.75'
$ws.Range('E22').Value = '  +7.64%  '
$ws.Range('D23').Value = '4.37'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('D24').Value = '2.12'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('D25').Value = '149.27'
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').Value = '6.87'
$ws.Range('E26').Value = '  -1.16%  '
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  -0.78%  '
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('D33').Value = '1.465.28'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').Value = '3.09'
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').Value = '0.876'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.558'
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0167'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = '0.915'
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').Value = '1.01'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').Value = '67.89'
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('E44').Value = '  +1.79%  '
$ws.Range('D45').Value = '2.47'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('D46').Value = '5.37'
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('D47').Value = '1.771.76'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('D49').Value = '87.31'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').Value = '0.0994'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₇0995'
$ws.Range('E51').Value = '  -6.12%  '

# The original cells used the default style (no custom format), so restore
# it now that the text values are safely stored -- keeps formatting identical.
$ws.Range("D2:E51").Style = "Normal"

